# Oppgave 1 og 2 fullført, mangler å oppdatere antall ledige plasser
#
# Adds a new "svar" column header to the "soknad" sheet (N1), copying the
# bold/bordered header style from the preceding header cell (M1), and
# updates the active selection like the author left it (P2).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing header's formatting onto the new header cell, then
# write its text. Copy/PasteSpecial (paste formats only) carries over the
# bold font + border + alignment style used by the rest of row 1.
$ws.Range("M1").Copy()
$ws.Range("N1").PasteSpecial(-4122)
$ws.Range("N1").Value = "svar"

# Leave the selection where the author ended up.
$ws.Range("P2").Select() | Out-Null
